# ----------------------------------------------------------------------
# Sprint plan edits for projectmanagement.docx
# ----------------------------------------------------------------------

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Narrative paragraph edits (also merges split runs into one run)
# ------------------------------------------------------------------

$old1 = "Students has the rights to know what does a University plan to do regarding on its developments and projects, but students only express their thoughts and ask questions on what does a project benefits for them and how can that project help the University by posting in different kinds of fan pages/ school pages that is related to the University. They often ask to the University government and make some protest for the purpose of it because they are wondering what is it for, why does that certain project costs that much and why did the University make that certain project. With these types of concerns, by the help of this Online Inquiry, they can post their questions and concerns about the University with the use of internet. In every questions that they posted, admin from that University can answer their questions and for them to know immediately what they ask for in a short amount of time."
$new1 = "Students have the right to know what their University plans to do regarding on its developments and projects, but students only express their thoughts and ask questions on what does a project benefits for them and how can that project help the University by posting in different kinds of fan pages/ school pages that is related to the University. They often ask to the University government and make some protest for the purpose of it because they are wondering what is it for, why does that certain project costs that much and why did the University make that certain project. With these types of concerns, by the help of this Online Inquiry, they can post their questions and concerns about the University with the use of internet. In every questions that they posted, admin from that University can answer their questions and for them to know immediately what they ask for in a short amount of time."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

$old2 = "The system is in public, so there will be no need for the students to create an account just to post a question.  "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2)

$old3 = "The system is an FAQ module of the website and users can post as many questions as they want. But administrators can only answer once in every post."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2)

$old4 = "There will be a category for the questions to sort and view them by the selected category of it."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $old4, 2)

# ------------------------------------------------------------------
# 2. Sprint-plan table (Table 1): update dates / durations per row,
#    then drop the trailing "Status" column and shrink the table width.
# ------------------------------------------------------------------

$sprintTable = $d.Tables.Item(1)

# NOTE: wdReplaceAll (Replace:=2) ignores the boundaries of the Range it
# is invoked on and replaces every match in the whole story/document, so
# for these short, frequently-duplicated cell values (e.g. "1/30") we
# scope the search to an explicit Range and use wdReplaceOne (Replace:=1)
# which correctly confines the match to that Range.
function Replace-CellText($table, $row, $col, $oldText, $newText) {
    $cellRange = $table.Cell($row, $col).Range
    $rng = $d.Range($cellRange.Start, $cellRange.End)
    $result = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)
    Write-Output "cell($row,$col) '$oldText'->'$newText': $result"
}

# Conceptualization of database model
Replace-CellText $sprintTable 3 3 "1/30" "2/9"
Replace-CellText $sprintTable 3 4 "2/1" "2/9"
Replace-CellText $sprintTable 3 5 "3 days" "1-2 hour/s"

# Client-side GUI
Replace-CellText $sprintTable 4 3 "1/30" "2/9"
Replace-CellText $sprintTable 4 4 "2/2" "2/12"
Replace-CellText $sprintTable 4 5 "4 days" "3 days"

# Admin-side GUI
Replace-CellText $sprintTable 5 3 "1/30" "2/9"
Replace-CellText $sprintTable 5 4 "2/2" "2/12"
Replace-CellText $sprintTable 5 5 "4 days" "3 days"

# CRUD Admin
Replace-CellText $sprintTable 7 3 "2/1" "2/13"
Replace-CellText $sprintTable 7 4 "2/5" "2/17"
Replace-CellText $sprintTable 7 5 "5 days" "4 days"

# Login/Logout
Replace-CellText $sprintTable 8 3 "2/6" "2/18"
Replace-CellText $sprintTable 8 4 "2/7" "2/19"

# CRUD Inquiry
Replace-CellText $sprintTable 9 3 "2/1" "2/13"
Replace-CellText $sprintTable 9 4 "2/5" "2/17"
Replace-CellText $sprintTable 9 5 "5 days" "4 days"

# Testing
Replace-CellText $sprintTable 11 3 "2/15" "2/20"
Replace-CellText $sprintTable 11 4 "2/17" "2/23"
Replace-CellText $sprintTable 11 5 "2 days" "3 days"

# Drop the "Status" column (6th column) entirely and shrink table width
$sprintTable.Columns.Item(6).Delete()
$sprintTable.PreferredWidth = 419.25

# ------------------------------------------------------------------
# 3. Difficulty / Estimation table (Table 2)
# ------------------------------------------------------------------

$diffTable = $d.Tables.Item(2)

Replace-CellText $diffTable 2 3 "2-3 days" "1-2 hours"
Replace-CellText $diffTable 3 3 "3-4 days" "3-4 days"
Replace-CellText $diffTable 4 3 "1 day" "1 day"
Replace-CellText $diffTable 5 2 "Small" "Medium"
Replace-CellText $diffTable 5 3 "4-5 days" "3-4 days"
Replace-CellText $diffTable 6 3 "4-5 days" "3-4 days"
Replace-CellText $diffTable 7 3 "1-2 days" "2-3 days"

# ------------------------------------------------------------------
# 4. Trailing empty paragraph added just before the final section break
# ------------------------------------------------------------------

$extraPara = $d.Paragraphs.Add()
$extraPara.Style = "Normal"
